$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRU")

$data = @(
    @(24, "tanh", 3, 0.0001969071360761338, [double]"3.295071680202894e-07", 0.44, 0.02085443213582039, 0.004968228749930859, 0.2349011407684074, 0.01640518580824816),
    @(34, "tanh", 4, [double]"4.954308708654279e-05", [double]"1.027730287589487e-06", 0.46, 0.01676929742097855, 0.003693493083119392, 0.2419355726006238, 0.007279144241870948),
    @(35, "tanh", 4, [double]"4.63090461339639e-05", [double]"1.612349719343376e-05", 0.08, 0.005075340159237385, 0.003563524456694722, 0.244540215561469, 0.0698721133214506),
    @(45, "tanh", 4, [double]"1.081165724733026e-06", [double]"1.368856402116335e-06", 0.15, 0.005473521072417498, 0.00330555415712297, 0.2389582113883663, 0.004937563463144684),
    @(45, "tanh", 4, [double]"8.632214722101613e-07", [double]"1.213804620863432e-06", 0.14, 0.005315640475600958, 0.003446863498538733, 0.2488232940825086, 0.02959955547877408),
    @(46, "tanh", 4, [double]"1.047117208703457e-05", [double]"1.061045628337706e-05", 0.25, 0.007363975513726473, 0.003657443448901176, 0.2395274915339935, 0.02731770837599144),
    @(36, "tanh", 4, [double]"3.35233887351799e-06", [double]"1.72935903209589e-06", 0.15, 0.00599592924118042, 0.00379849155433476, 0.2415428440029574, 0.3740611493409658),
    @(42, "tanh", 4, [double]"3.278416300792034e-07", [double]"3.903801849116737e-08", 0.22, 0.006886400282382965, 0.003267077961936593, 0.2417079557657851, 0.007579290575687332),
    @(29, "tanh", 4, [double]"1.034125341900445e-07", [double]"3.939985410670502e-08", 0.1, 0.005337257869541645, 0.003191662020981312, 0.2454559858494016, 0.1035385409843182),
    @(43, "tanh", 4, [double]"3.139947723386321e-07", [double]"2.519928634017416e-07", 0.22, 0.006793499458581209, 0.003268587635830045, 0.2410615656093686, 0.007396387556035302)
)

$startRow = 29
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$wb.Save()
